$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new activity log entry (row 10): Di, 26.11.2019, 07:45-12:30,
# "Benutzer und Kinder auflisten"
$ws.Range("A10").Value = "Di"
$ws.Range("B10").Value = "26.11.2019"
$ws.Range("C10").Value = 0.32291666666666669
$ws.Range("D10").Value = 0.52083333333333337
$ws.Range("E10").Value = "Benutzer und Kinder auflisten"

# Move the active selection to A11, as in the saved workbook
$ws.Range("A11").Select()
